$d = $word.ActiveDocument

function Find-ParagraphByText($needle) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like $needle) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) PEC_COMANDO mergefield -> hardcoded PEC address (plain italic text run)
# ---------------------------------------------------------------------------
$p1 = Find-ParagraphByText("*PEC_COMANDO*")
$r1 = $d.Range($p1.Range.Start, $p1.Range.End)
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
    '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:widowControl w:val="false"/><w:bidi w:val="0"/><w:jc w:val="center"/><w:rPr><w:i/><w:i/><w:iCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>com.savona@cert.vigilfuoco.it</w:t></w:r>' + `
    '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$r1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 2) "RGNR n°..." -> "RGNR n°" + MERGEFIELD $RGNR
# ---------------------------------------------------------------------------
$p2 = Find-ParagraphByText("*RGNR*")
$r2 = $d.Range($p2.Range.Start, $p2.Range.End)
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
    '<w:p><w:pPr><w:pStyle w:val="Contenutotabella"/><w:widowControl w:val="false"/><w:bidi w:val="0"/><w:jc w:val="center"/><w:rPr/></w:pPr>' + `
    '<w:r><w:rPr/><w:t>RGNR n&#176;</w:t></w:r>' + `
    '<w:r><w:rPr/><w:fldChar w:fldCharType="begin"/></w:r>' + `
    '<w:r><w:rPr/><w:instrText xml:space="preserve"> MERGEFIELD $RGNR </w:instrText></w:r>' + `
    '<w:r><w:rPr/><w:fldChar w:fldCharType="separate"/></w:r>' + `
    '<w:r><w:rPr/><w:t>&lt;$RGNR&gt;</w:t></w:r>' + `
    '<w:r><w:rPr/><w:fldChar w:fldCharType="end"/></w:r>' + `
    '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$r2.InsertXML($xml2)

# ---------------------------------------------------------------------------
# 3) Add a new paragraph "(rif. PM <$MAGISTRATO>)" after the "...Tribunale di
#    <$TRIBUNALE>" paragraph, inside the same table cell.
# ---------------------------------------------------------------------------
$p3 = Find-ParagraphByText("*Tribunale*")
$r3 = $d.Range($p3.Range.End, $p3.Range.End)
$xml3 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
    '<w:p><w:pPr><w:pStyle w:val="Contenutotabella"/><w:widowControl w:val="false"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr>' + `
    '<w:r><w:rPr/><w:t xml:space="preserve">(rif. PM </w:t></w:r>' + `
    '<w:r><w:rPr/><w:fldChar w:fldCharType="begin"/></w:r>' + `
    '<w:r><w:rPr/><w:instrText xml:space="preserve"> MERGEFIELD $MAGISTRATO </w:instrText></w:r>' + `
    '<w:r><w:rPr/><w:fldChar w:fldCharType="separate"/></w:r>' + `
    '<w:r><w:rPr/><w:t>&lt;$MAGISTRATO&gt;</w:t></w:r>' + `
    '<w:r><w:rPr/><w:fldChar w:fldCharType="end"/></w:r>' + `
    '<w:r><w:rPr/><w:t>)</w:t></w:r>' + `
    '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$r3.InsertXML($xml3)

Write-Output "Edits applied."
